$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 144, duplicating the existing row 144
# (this shifts the former row 144 and all following rows down by one).
$ws.Rows("144:144").Copy()
$ws.Rows("144:144").Insert()

# Now update the newly inserted row 144 with the new data values.
$ws.Range("D144").Value = 44566
$ws.Range("J144").Value = 200
$ws.Range("K144").Value = 4000
$ws.Range("L144").Value = 4500
$ws.Range("M144").Value = 4300
$ws.Range("P144").Value = 1433
